# ENKAVI_VariableInfo.xlsx update
# - insert a new "item_num" column just before the "comment" column on every
#   sheet that has that column (Measures, ID, Dems, Dates, NewVars)
# - populate item_num values on the Measures sheet
# - drop the now-unused trailing blank rows on Measures
# - clear the leftover number-format / font styling that was only used by
#   those blank rows / cells
# - turn off the AutoFilter on Measures (range no longer matches filter db)
# - widen the Measures!_xlnm._FilterDatabase defined name to include the
#   new column
# - tidy up a couple of view/selection details

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Measures sheet
# ---------------------------------------------------------------------------
$measures = $wb.Worksheets.Item("Measures")

# Insert the new column (O) ahead of "comment" (was O, becomes P)
$measures.Columns.Item(15).Insert()
$measures.Range("O1").Value = "item_num"

# Fill in item_num for every data row (grouped by identical value)
$measures.Range("O2:O23").Value = 1
$measures.Range("O24:O25").Value = 6
$measures.Range("O26:O29").Value = 30
$measures.Range("O30:O33").Value = 32
$measures.Range("O34:O35").Value = 10
$measures.Range("O36:O39").Value = 6
$measures.Range("O40:O41").Value = 1
$measures.Range("O42:O43").Value = 6

# The old trailing blank/formatted rows (44-57) are no longer needed
$measures.Range("A44:A57").EntireRow.Delete()

# Clear the leftover styling (applyFont / applyNumberFormat) that used to
# live on column J and M26:M35.  Cells that still hold a value only need
# their formatting cleared; cells that were empty (style-only) go away
# entirely once both value and format are gone.
$measures.Range("J2:J23").ClearFormats()
$measures.Range("J40:J41").ClearFormats()
$measures.Range("J24:J25").Clear()
$measures.Range("J36:J39").Clear()
$measures.Range("J42:J43").Clear()
$measures.Range("M26:M35").ClearFormats()

# The data no longer needs to be auto-filtered
$measures.AutoFilterMode = $false

# Keep the filter-database defined name in sync with the extra column
$wb.Names.Item("Measures!_FilterDatabase").RefersTo = "=Measures!`$A`$1:`$P`$43"

$measures.Activate()
$measures.Range("K16").Select()

# ---------------------------------------------------------------------------
# ID sheet
# ---------------------------------------------------------------------------
$id = $wb.Worksheets.Item("ID")
$id.Columns.Item(15).Insert()
$id.Range("O1").Value = "item_num"
$id.Activate()
$id.Range("L15").Select()

# ---------------------------------------------------------------------------
# Dems sheet
# ---------------------------------------------------------------------------
$dems = $wb.Worksheets.Item("Dems")
$dems.Columns.Item(15).Insert()
$dems.Range("O1").Value = "item_num"
$dems.Activate()
$dems.Range("M10").Select()

# ---------------------------------------------------------------------------
# Dates sheet
# ---------------------------------------------------------------------------
$dates = $wb.Worksheets.Item("Dates")
$dates.Columns.Item(15).Insert()
$dates.Range("O1").Value = "item_num"
$dates.Activate()
$dates.Range("N6").Select()

# ---------------------------------------------------------------------------
# NewVars sheet (has an extra "varfolder" column, so item_num lands in P)
# ---------------------------------------------------------------------------
$newVars = $wb.Worksheets.Item("NewVars")
$newVars.Columns.Item(16).Insert()
$newVars.Range("P1").Value = "item_num"
$newVars.Activate()
$newVars.Range("O9").Select()

# Re-activate Measures as the originally active sheet
$measures.Activate()
